$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = 9948.2999999999993
$ws.Range("B6").Value = 10049.799999999999
$ws.Range("C6").Value = 286
$ws.Range("D6").Value = 283.11
$ws.Range("E6").Value = $false
$ws.Range("F6").Value = -1.01
$ws.Range("G6").Value = 42612.675081018519
$ws.Range("G6").NumberFormat = "m/d/yy h:mm"
$ws.Range("H6").Value = $false
